# Normalize the "Валута" (currency) column: every data row in column C
# (rows 2-23) is set to "ДЕН", replacing whichever currency ("ЕУР"/"УСД")
# previously occupied that cell. Once nothing references "ЕУР"/"УСД" any
# more, those now-unused shared strings are dropped automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 23; $r++) {
    $ws.Range("C$r").Value = "ДЕН"
}

# Match the author's last on-screen selection recorded in the saved file.
$ws.Range("D31").Select()
